$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the hard-coded offset value in I8 (was 6, should be 5).
# This also recalculates the dependent formulas in I5 and I6.
$ws.Range("I8").Value = 5

# Recalculate so the shared/regular formulas pick up the new value.
$excel.Calculate()

# Update the current selection to match the saved view state (I9).
$ws.Activate()
$ws.Range("I9").Select()
